# The deck's single theme ("Integral", ppt/theme/theme1.xml) is being
# re-coloured to the stock "Office Theme" palette (whose colours
# previously lived in ppt/theme/theme2.xml, used only by the notes
# master). The font scheme and format scheme are already identical
# between the two themes, so only the 12 theme colours differ.
#
# Helper: convert an "RRGGBB" hex string into the BGR-packed decimal
# long that PowerPoint's ColorFormat.RGB property expects.
function Hex2BGR($hex) {
    $rr = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $gg = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $bb = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return ($bb * 65536) + ($gg * 256) + $rr
}

$p = $ppt.ActivePresentation

# Any slide exposes the one shared ThemeColorScheme (12 theme colours,
# in dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink order).
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$officeColors = @(
    "000000",  # 1  dk1
    "FFFFFF",  # 2  lt1
    "44546A",  # 3  dk2
    "E7E6E6",  # 4  lt2
    "5B9BD5",  # 5  accent1
    "ED7D31",  # 6  accent2
    "A5A5A5",  # 7  accent3
    "FFC000",  # 8  accent4
    "4472C4",  # 9  accent5
    "70AD47",  # 10 accent6
    "0563C1",  # 11 hlink
    "954F72"   # 12 folHlink
)

for ($i = 0; $i -lt $officeColors.Count; $i++) {
    $tcs.Item($i + 1).RGB = Hex2BGR $officeColors[$i]
}
